$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "Study_Types" (C) / "QA_Excel_Files" (D) columns are dropped
# entirely. Deleting them with a shift-left brings the
# "manage_population_file_name" / "manage_population_file_to_upload"
# columns (previously E:F) into C:D, and the old
# "Expected_File_names" / "Files_to_upload" columns (previously G:H)
# into E:F - exactly matching the new A1:F4 layout.
$ws.Range("C1:D5").Delete(-4163)   # xlShiftToLeft

# The full local "D:\VersionControl\..." paths are replaced with the
# short relative "\ExtractionTemplate\..." form.
$ws.Range("D2").Value = "\ExtractionTemplate\sample_population_data_sheet.xlsx"
$ws.Range("F2").Value = "\ExtractionTemplate\sample_population_data_sheet.xlsx"
$ws.Range("F3").Value = "\ExtractionTemplate\Pfizer-Economic-Report-20220725120854.xlsx"

# Header row: centered, no wrap.
$ws.Range("A1:F1").WrapText = $false
$ws.Range("A1:F1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:F1").VerticalAlignment = -4108     # xlCenter

# Data rows: default alignment, no wrap. Only touch the cells that are
# actually populated so no stray, empty-but-styled cells get written out
# for the gaps (C3:D3, C4:F4).
$ws.Range("A2:F2,A3:B3,E3:F3,A4:B4").WrapText = $false

# Row 5 (which only ever held the old Study_Types/QA_Excel_Files overflow)
# is now completely empty - drop it so the used range ends at row 4.
$ws.Rows("5").Delete()

# The rows no longer need the tall, wrap-driven 57.6pt height - let Excel
# recompute the natural (default) row height for the new content.
$ws.Rows("1:4").AutoFit()

# Re-fit the columns to the new, shorter content.
$ws.Columns("A:F").AutoFit()

$ws.Range("D11").Select()
